$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "查宇航"
$ws.Range("B2").Value2 = 201917
$ws.Range("C2").Value2 = 97.5
$ws.Range("D2").Value2 = 97
$ws.Range("E2").Value2 = 96

$ws.Range("A3").Value2 = "钱梓暄"
$ws.Range("B3").Value2 = 201932
$ws.Range("C3").Value2 = 96
$ws.Range("D3").Value2 = 97
$ws.Range("E3").Value2 = 96

$ws.Range("A4").Value2 = "王慧茹"
$ws.Range("B4").Value2 = 201923
$ws.Range("C4").Value2 = 98
$ws.Range("D4").Value2 = 99
$ws.Range("E4").Value2 = 91

$ws.Range("A5").Value2 = "刘可馨"
$ws.Range("B5").Value2 = 201927
$ws.Range("C5").Value2 = 95
$ws.Range("D5").Value2 = 100
$ws.Range("E5").Value2 = 90

$ws.Range("A6").Value2 = "李文琪"
$ws.Range("B6").Value2 = 201912
$ws.Range("C6").Value2 = 97
$ws.Range("D6").Value2 = 91
$ws.Range("E6").Value2 = 96

$ws.Range("A7").Value2 = "陶永胜"
$ws.Range("B7").Value2 = 201916
$ws.Range("C7").Value2 = 98
$ws.Range("D7").Value2 = 95
$ws.Range("E7").Value2 = 90

$ws.Range("A8").Value2 = "毛瑜彤"
$ws.Range("B8").Value2 = 201925
$ws.Range("C8").Value2 = 96
$ws.Range("D8").Value2 = 96
$ws.Range("E8").Value2 = 91

$ws.Range("A9").Value2 = "蔡奕扬"
$ws.Range("B9").Value2 = 201902
$ws.Range("C9").Value2 = 89
$ws.Range("D9").Value2 = 96
$ws.Range("E9").Value2 = 98

$ws.Range("A10").Value2 = "屠佳佳"
$ws.Range("B10").Value2 = 201934
$ws.Range("C10").Value2 = 92.5
$ws.Range("D10").Value2 = 92
$ws.Range("E10").Value2 = 96

$ws.Range("A11").Value2 = "王承彦"
$ws.Range("B11").Value2 = 201910
$ws.Range("C11").Value2 = 94.5
$ws.Range("D11").Value2 = 90
$ws.Range("E11").Value2 = 96

$ws.Range("A12").Value2 = "霍懿昕"
$ws.Range("B12").Value2 = 201926
$ws.Range("C12").Value2 = 97.5
$ws.Range("D12").Value2 = 86
$ws.Range("E12").Value2 = 95

$ws.Range("A13").Value2 = "王斌臣"
$ws.Range("B13").Value2 = 201903
$ws.Range("C13").Value2 = 93.5
$ws.Range("D13").Value2 = 94
$ws.Range("E13").Value2 = 90

$ws.Range("A14").Value2 = "朱国煜"
$ws.Range("B14").Value2 = 201906
$ws.Range("C14").Value2 = 94
$ws.Range("D14").Value2 = 97
$ws.Range("E14").Value2 = 86

$ws.Range("A15").Value2 = "华奕轩"
$ws.Range("B15").Value2 = 201921
$ws.Range("C15").Value2 = 97.5
$ws.Range("D15").Value2 = 77
$ws.Range("E15").Value2 = 100

$ws.Range("A16").Value2 = "施浩林"
$ws.Range("B16").Value2 = 201918
$ws.Range("C16").Value2 = 90
$ws.Range("D16").Value2 = 95
$ws.Range("E16").Value2 = 89

$ws.Range("A17").Value2 = "周越芊"
$ws.Range("B17").Value2 = 201931
$ws.Range("C17").Value2 = 92.5
$ws.Range("D17").Value2 = 90
$ws.Range("E17").Value2 = 90

$ws.Range("A18").Value2 = "包静宜"
$ws.Range("B18").Value2 = 201928
$ws.Range("C18").Value2 = 97
$ws.Range("D18").Value2 = 79
$ws.Range("E18").Value2 = 94

$ws.Range("A19").Value2 = "葛逸玮"
$ws.Range("B19").Value2 = 201920
$ws.Range("C19").Value2 = 94.5
$ws.Range("D19").Value2 = 77
$ws.Range("E19").Value2 = 94

$ws.Range("A20").Value2 = "毛骐耀"
$ws.Range("B20").Value2 = 201905
$ws.Range("C20").Value2 = 88
$ws.Range("D20").Value2 = 97
$ws.Range("E20").Value2 = 79

$ws.Range("A21").Value2 = "金佳琪"
$ws.Range("B21").Value2 = 201937
$ws.Range("C21").Value2 = 86
$ws.Range("D21").Value2 = 88
$ws.Range("E21").Value2 = 90

$ws.Range("A22").Value2 = "蒋梓昊"
$ws.Range("B22").Value2 = 201907
$ws.Range("C22").Value2 = 94.5
$ws.Range("D22").Value2 = 97
$ws.Range("E22").Value2 = 72

$ws.Range("A23").Value2 = "张宇凡"
$ws.Range("B23").Value2 = 201933
$ws.Range("C23").Value2 = 89.5
$ws.Range("D23").Value2 = 87
$ws.Range("E23").Value2 = 86

$ws.Range("A24").Value2 = "罗倩霞"
$ws.Range("B24").Value2 = 201924
$ws.Range("C24").Value2 = 91
$ws.Range("D24").Value2 = 87
$ws.Range("E24").Value2 = 83

$ws.Range("A25").Value2 = "杨孟凡"
$ws.Range("B25").Value2 = 201914
$ws.Range("C25").Value2 = 92
$ws.Range("D25").Value2 = 69
$ws.Range("E25").Value2 = 96

$ws.Range("A26").Value2 = "沈一凡"
$ws.Range("B26").Value2 = 201911
$ws.Range("C26").Value2 = 92.5
$ws.Range("D26").Value2 = 95
$ws.Range("E26").Value2 = 69

$ws.Range("A27").Value2 = "潘音琪"
$ws.Range("B27").Value2 = 201936
$ws.Range("C27").Value2 = 88.5
$ws.Range("D27").Value2 = 79
$ws.Range("E27").Value2 = 87

$ws.Range("A28").Value2 = "杨之文"
$ws.Range("B28").Value2 = 201929
$ws.Range("C28").Value2 = 68
$ws.Range("D28").Value2 = 93
$ws.Range("E28").Value2 = 93

$ws.Range("A29").Value2 = "邹锐"
$ws.Range("B29").Value2 = 201919
$ws.Range("C29").Value2 = 86.5
$ws.Range("D29").Value2 = 81
$ws.Range("E29").Value2 = 86

$ws.Range("A30").Value2 = "席庆"
$ws.Range("B30").Value2 = 201908
$ws.Range("C30").Value2 = 83
$ws.Range("D30").Value2 = 77
$ws.Range("E30").Value2 = 92

$ws.Range("A31").Value2 = "金熙晨"
$ws.Range("B31").Value2 = 201915
$ws.Range("C31").Value2 = 89
$ws.Range("D31").Value2 = 76
$ws.Range("E31").Value2 = 84

$ws.Range("A32").Value2 = "戴逸辰"
$ws.Range("B32").Value2 = 201922
$ws.Range("C32").Value2 = 65
$ws.Range("D32").Value2 = 93
$ws.Range("E32").Value2 = 91

$ws.Range("A33").Value2 = "周智宸"
$ws.Range("B33").Value2 = 201909
$ws.Range("C33").Value2 = 85
$ws.Range("D33").Value2 = 68
$ws.Range("E33").Value2 = 92

$ws.Range("A34").Value2 = "颜宇晨"
$ws.Range("B34").Value2 = 201904
$ws.Range("C34").Value2 = 78.5
$ws.Range("D34").Value2 = 79
$ws.Range("E34").Value2 = 84

$ws.Range("A35").Value2 = "张晨阳"
$ws.Range("B35").Value2 = 201901
$ws.Range("C35").Value2 = 92
$ws.Range("D35").Value2 = 68
$ws.Range("E35").Value2 = 78

$ws.Range("A36").Value2 = "陈沫兰"
$ws.Range("B36").Value2 = 201930
$ws.Range("C36").Value2 = 87.5
$ws.Range("D36").Value2 = 61
$ws.Range("E36").Value2 = 68

$ws.Range("A37").Value2 = "陆俊祺"
$ws.Range("B37").Value2 = 201913
$ws.Range("C37").Value2 = 83
$ws.Range("D37").Value2 = 55
$ws.Range("E37").Value2 = 64

$ws.Range("A38").Value2 = "程梓涵"
$ws.Range("B38").Value2 = 201938
$ws.Range("C38").Value2 = -1
$ws.Range("D38").Value2 = -1
$ws.Range("E38").Value2 = -1

$ws.Range("A39").Value2 = "沈琪"
$ws.Range("B39").Value2 = 201935
$ws.Range("C39").Value2 = -1
$ws.Range("D39").Value2 = -1
$ws.Range("E39").Value2 = -1
